$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1448.4546
$ws.Range("I28").Value = 1554.3125
$ws.Range("J28").Value = 1166.1666
$ws.Range("K28").Value = 1554.3125
$ws.Range("L28").Value = 1166.1666
$ws.Range("M28").Value = -1069.3125
$ws.Range("N28").Value = -2136.1666
$ws.Range("H29").Value = 4555.25
$ws.Range("I29").Value = 1110.5
$ws.Range("K29").Value = 3331.5
$ws.Range("M29").Value = -3050.5
$ws.Range("H38").Value = 6576.385
$ws.Range("J38").Value = 10298.625
$ws.Range("L38").Value = 30895.875
$ws.Range("N38").Value = -31639.875
$ws.Range("H58").Value = 3943.7273
$ws.Range("I58").Value = 642
$ws.Range("J58").Value = 9721.75
$ws.Range("K58").Value = 1926
$ws.Range("L58").Value = 29165.25
$ws.Range("M58").Value = -1776
$ws.Range("N58").Value = -29465.25
$ws.Range("H69").Value = 2011200.6
$ws.Range("I69").Value = 2011200.6
$ws.Range("K69").Value = 6033601.800000001
$ws.Range("M69").Value = -6032727.800000001
$ws.Range("H72").Value = 2011200.6
$ws.Range("I72").Value = 2011200.6
$ws.Range("K72").Value = 18100805.4
$ws.Range("M72").Value = -18096437.4
$ws.Range("H92").Value = 196183.75
$ws.Range("J92").Value = 350195.75
$ws.Range("L92").Value = 350195.75
$ws.Range("N92").Value = -352691.75
$ws.Range("H94").Value = 2023.6
$ws.Range("J94").Value = 3424.5
$ws.Range("L94").Value = 3424.5
$ws.Range("N94").Value = -4326.5
$ws.Range("H96").Value = 518.1818
$ws.Range("I96").Value = 477.15384
$ws.Range("J96").Value = 577.44446
$ws.Range("K96").Value = 1431.46152
$ws.Range("L96").Value = 1732.33338
$ws.Range("M96").Value = -58.46152000000006
$ws.Range("N96").Value = -4478.33338
$ws.Range("H103").Value = 1004.1905
$ws.Range("I103").Value = 583.25
$ws.Range("J103").Value = 1103.2354
$ws.Range("K103").Value = 1749.75
$ws.Range("L103").Value = 3309.7062
$ws.Range("M103").Value = -1163.75
$ws.Range("N103").Value = -4481.706200000001
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H107").Value = 786
$ws.Range("I107").Value = 489.86957
$ws.Range("K107").Value = 489.86957
$ws.Range("M107").Value = 1430.13043
$ws.Range("H112").Value = 25206.158
$ws.Range("J112").Value = 29713.812
$ws.Range("L112").Value = 89141.436
$ws.Range("N112").Value = -91357.436
$ws.Range("H116").Value = 6354.115
$ws.Range("I116").Value = 5420.0625
$ws.Range("K116").Value = 5420.0625
$ws.Range("M116").Value = -1978.0625
$ws.Range("H132").Value = 1798.7291
$ws.Range("I132").Value = 1589.55
$ws.Range("J132").Value = 2844.625
$ws.Range("K132").Value = 4768.65
$ws.Range("L132").Value = 8533.875
$ws.Range("M132").Value = -2238.65
$ws.Range("N132").Value = -13593.875
$ws.Range("H135").Value = 1185.5714
$ws.Range("I135").Value = 1121.4828
$ws.Range("J135").Value = 1495.3334
$ws.Range("K135").Value = 10093.3452
$ws.Range("L135").Value = 13458.0006
$ws.Range("M135").Value = -7558.3452
$ws.Range("N135").Value = -18528.0006
$ws.Range("H138").Value = 2345.5847
$ws.Range("I138").Value = 765.9697
$ws.Range("J138").Value = 3974.5625
$ws.Range("K138").Value = 2297.9091
$ws.Range("L138").Value = 11923.6875
$ws.Range("M138").Value = 2842.0909
$ws.Range("N138").Value = -22203.6875
$ws.Range("H139").Value = 70699.38
$ws.Range("J139").Value = 70699.38
$ws.Range("L139").Value = 70699.38
$ws.Range("N139").Value = -80979.38
$ws.Range("H140").Value = 88966.2
$ws.Range("J140").Value = 88966.2
$ws.Range("L140").Value = 88966.2
$ws.Range("N140").Value = -99326.2
$ws.Range("H141").Value = 4516
$ws.Range("I141").Value = 4379.647
$ws.Range("K141").Value = 13138.941
$ws.Range("M141").Value = -7958.940999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1074.6389
$ws.Range("I2").Value = 1193.5
$ws.Range("J2").Value = 480.33334
$ws.Range("K2").Value = 1193.5
$ws.Range("L2").Value = 480.33334
$ws.Range("M2").Value = -1080.5
$ws.Range("N2").Value = -706.33334
$ws.Range("H32").Value = 12146.553
$ws.Range("I32").Value = 6213.147
$ws.Range("J32").Value = 62580.5
$ws.Range("K32").Value = 6213.147
$ws.Range("L32").Value = 62580.5
$ws.Range("M32").Value = -5926.147
$ws.Range("N32").Value = -63154.5
$ws.Range("H45").Value = 2046.4706
$ws.Range("I45").Value = 1342.2858
$ws.Range("K45").Value = 1342.2858
$ws.Range("M45").Value = -965.2858000000001
$ws.Range("H61").Value = 28055.162
$ws.Range("I61").Value = 2011.091
$ws.Range("K61").Value = 2011.091
$ws.Range("M61").Value = -1799.091
$ws.Range("H74").Value = 80479.39
$ws.Range("I74").Value = 58063.562
$ws.Range("J74").Value = 131715.58
$ws.Range("K74").Value = 58063.562
$ws.Range("L74").Value = 131715.58
$ws.Range("M74").Value = -57189.562
$ws.Range("N74").Value = -133463.58
$ws.Range("H77").Value = 80479.39
$ws.Range("I77").Value = 58063.562
$ws.Range("J77").Value = 131715.58
$ws.Range("K77").Value = 290317.81
$ws.Range("L77").Value = 658577.8999999999
$ws.Range("M77").Value = -285949.81
$ws.Range("N77").Value = -667313.8999999999
$ws.Range("H80").Value = 34056.2
$ws.Range("J80").Value = 40070.25
$ws.Range("L80").Value = 40070.25
$ws.Range("N80").Value = -42066.25
$ws.Range("H83").Value = 34056.2
$ws.Range("J83").Value = 40070.25
$ws.Range("L83").Value = 120210.75
$ws.Range("N83").Value = -130194.75
$ws.Range("H97").Value = 47619964
$ws.Range("I97").Value = 949.1429
$ws.Range("K97").Value = 949.1429
$ws.Range("M97").Value = -453.1429000000001
$ws.Range("H102").Value = 26411634
$ws.Range("I102").Value = 34334696
$ws.Range("K102").Value = 34334696
$ws.Range("M102").Value = -34333074
$ws.Range("H116").Value = 1074.6389
$ws.Range("I116").Value = 1193.5
$ws.Range("J116").Value = 480.33334
$ws.Range("K116").Value = 1193.5
$ws.Range("L116").Value = 480.33334
$ws.Range("M116").Value = 1100.5
$ws.Range("N116").Value = -5068.33334
$ws.Range("H122").Value = 101215.78
$ws.Range("I122").Value = 948.4
$ws.Range("K122").Value = 2845.2
$ws.Range("M122").Value = -395.1999999999998
$ws.Range("H132").Value = 1594.0571
$ws.Range("I132").Value = 1773.2069
$ws.Range("J132").Value = 728.1667
$ws.Range("K132").Value = 5319.620699999999
$ws.Range("L132").Value = 2184.5001
$ws.Range("M132").Value = -2789.620699999999
$ws.Range("N132").Value = -7244.5001
$ws.Range("H136").Value = 28055.162
$ws.Range("I136").Value = 2011.091
$ws.Range("K136").Value = 6033.272999999999
$ws.Range("M136").Value = -3483.272999999999
$ws.Range("H140").Value = 66272.5
$ws.Range("J140").Value = 74900
$ws.Range("L140").Value = 74900
$ws.Range("N140").Value = -85260
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1074.6389
$ws.Range("I3").Value = 1193.5
$ws.Range("J3").Value = 480.33334
$ws.Range("K3").Value = 1193.5
$ws.Range("L3").Value = 480.33334
$ws.Range("M3").Value = -1079.5
$ws.Range("N3").Value = -708.33334
$ws.Range("H82").Value = 25127.9
$ws.Range("J82").Value = 40757.09
$ws.Range("L82").Value = 40757.09
$ws.Range("N82").Value = -41523.09
$ws.Range("H85").Value = 25127.9
$ws.Range("J85").Value = 40757.09
$ws.Range("L85").Value = 40757.09
$ws.Range("N85").Value = -43409.09
$ws.Range("H86").Value = 48732.42
$ws.Range("I86").Value = 1447.875
$ws.Range("J86").Value = 300916.66
$ws.Range("K86").Value = 1447.875
$ws.Range("L86").Value = 300916.66
$ws.Range("M86").Value = -324.875
$ws.Range("N86").Value = -303162.66
$ws.Range("H89").Value = 48732.42
$ws.Range("I89").Value = 1447.875
$ws.Range("J89").Value = 300916.66
$ws.Range("K89").Value = 7239.375
$ws.Range("L89").Value = 1504583.3
$ws.Range("M89").Value = -1623.375
$ws.Range("N89").Value = -1515815.3
$ws.Range("H105").Value = 3074.75
$ws.Range("I105").Value = 2100
$ws.Range("K105").Value = 2100
$ws.Range("M105").Value = -353
$ws.Range("H134").Value = 1392.8975
$ws.Range("I134").Value = 1397
$ws.Range("J134").Value = 1357
$ws.Range("K134").Value = 4191
$ws.Range("L134").Value = 4071
$ws.Range("M134").Value = -1656
$ws.Range("N134").Value = -9141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 647
$ws.Range("I14").Value = 647
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 647
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -477
$ws.Range("N14").ClearContents()
$ws.Range("H16").Value = 2096
$ws.Range("J16").Value = 3396.5
$ws.Range("L16").Value = 3396.5
$ws.Range("N16").Value = -3970.5
$ws.Range("H58").Value = 5697.147
$ws.Range("I58").Value = 5571.7827
$ws.Range("J58").Value = 5959.273
$ws.Range("K58").Value = 5571.7827
$ws.Range("L58").Value = 5959.273
$ws.Range("M58").Value = -5368.7827
$ws.Range("N58").Value = -6365.273
$ws.Range("H86").Value = 9464.8
$ws.Range("I86").Value = 5634.364
$ws.Range("K86").Value = 5634.364
$ws.Range("M86").Value = -4511.364
$ws.Range("H89").Value = 9464.8
$ws.Range("I89").Value = 5634.364
$ws.Range("K89").Value = 28171.82
$ws.Range("M89").Value = -22555.82
$ws.Range("H99").Value = 3646.276
$ws.Range("I99").Value = 3340.577
$ws.Range("J99").Value = 6295.6665
$ws.Range("K99").Value = 3340.577
$ws.Range("L99").Value = 6295.6665
$ws.Range("M99").Value = -1842.577
$ws.Range("N99").Value = -9291.6665
$ws.Range("H105").Value = 1422.6154
$ws.Range("I105").Value = 1153.1818
$ws.Range("J105").Value = 2904.5
$ws.Range("K105").Value = 1153.1818
$ws.Range("L105").Value = 2904.5
$ws.Range("M105").Value = 593.8181999999999
$ws.Range("N105").Value = -6398.5
$ws.Range("H113").Value = 2096
$ws.Range("J113").Value = 3396.5
$ws.Range("L113").Value = 3396.5
$ws.Range("N113").Value = -7736.5
$ws.Range("H122").Value = 1891.5625
$ws.Range("I122").Value = 1752.9
$ws.Range("J122").Value = 2122.6667
$ws.Range("K122").Value = 5258.700000000001
$ws.Range("L122").Value = 6368.000100000001
$ws.Range("M122").Value = -2808.700000000001
$ws.Range("N122").Value = -11268.0001
$ws.Range("H126").Value = 3646.276
$ws.Range("I126").Value = 3340.577
$ws.Range("J126").Value = 6295.6665
$ws.Range("K126").Value = 10021.731
$ws.Range("L126").Value = 18886.9995
$ws.Range("M126").Value = -7551.731
$ws.Range("N126").Value = -23826.9995
$ws.Range("H132").Value = 2689.353
$ws.Range("I132").Value = 2816.8262
$ws.Range("J132").Value = 2422.818
$ws.Range("K132").Value = 8450.4786
$ws.Range("L132").Value = 7268.454000000001
$ws.Range("M132").Value = -5920.4786
$ws.Range("N132").Value = -12328.454
$ws.Range("H136").Value = 5697.147
$ws.Range("I136").Value = 5571.7827
$ws.Range("J136").Value = 5959.273
$ws.Range("K136").Value = 16715.3481
$ws.Range("L136").Value = 17877.819
$ws.Range("M136").Value = -14165.3481
$ws.Range("N136").Value = -22977.819
$ws.Range("H141").Value = 163678.3
$ws.Range("J141").Value = 163678.3
$ws.Range("L141").Value = 163678.3
$ws.Range("N141").Value = -174038.3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25335.893
$ws.Range("I4").Value = 14669.884
$ws.Range("K4").Value = 44009.652
$ws.Range("M4").Value = -43897.652
$ws.Range("H7").Value = 1942.2858
$ws.Range("I7").Value = 1739.8
$ws.Range("J7").Value = 2448.5
$ws.Range("K7").Value = 5219.4
$ws.Range("L7").Value = 7345.5
$ws.Range("M7").Value = -5107.4
$ws.Range("N7").Value = -7569.5
$ws.Range("H14").Value = 1313.4814
$ws.Range("I14").Value = 1313.4814
$ws.Range("K14").Value = 3940.4442
$ws.Range("M14").Value = -3767.4442
$ws.Range("H34").Value = 2114.9
$ws.Range("I34").Value = 274.5
$ws.Range("J34").Value = 2575
$ws.Range("K34").Value = 823.5
$ws.Range("L34").Value = 7725
$ws.Range("M34").Value = -739.5
$ws.Range("N34").Value = -7893
$ws.Range("H39").Value = 3540.1667
$ws.Range("J39").Value = 4372.8125
$ws.Range("L39").Value = 13118.4375
$ws.Range("N39").Value = -13706.4375
$ws.Range("H55").Value = 3549.25
$ws.Range("I55").Value = 1565.8334
$ws.Range("J55").Value = 9499.5
$ws.Range("K55").Value = 4697.5002
$ws.Range("L55").Value = 28498.5
$ws.Range("M55").Value = -4520.5002
$ws.Range("N55").Value = -28852.5
$ws.Range("H92").Value = 1376.3334
$ws.Range("I92").Value = 1164.6666
$ws.Range("J92").Value = 1482.1666
$ws.Range("K92").Value = 3493.9998
$ws.Range("L92").Value = 4446.4998
$ws.Range("M92").Value = -2245.9998
$ws.Range("N92").Value = -6942.4998
$ws.Range("H121").Value = 7696813.5
$ws.Range("I121").Value = 13999.667
$ws.Range("J121").Value = 10001658
$ws.Range("K121").Value = 41999.001
$ws.Range("L121").Value = 30004974
$ws.Range("M121").Value = -40689.001
$ws.Range("N121").Value = -30007594
$ws.Range("H137").Value = 5063.383
$ws.Range("I137").Value = 6431.6665
$ws.Range("J137").Value = 4970.091
$ws.Range("K137").Value = 19294.9995
$ws.Range("L137").Value = 14910.273
$ws.Range("M137").Value = -14194.9995
$ws.Range("N137").Value = -25110.273
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19516.262
$ws.Range("I46").Value = 9306
$ws.Range("J46").Value = 23983.25
$ws.Range("K46").Value = 9306
$ws.Range("L46").Value = 23983.25
$ws.Range("M46").Value = -9150
$ws.Range("N46").Value = -24295.25
$ws.Range("H70").Value = 7156.9
$ws.Range("I70").Value = 5938.857
$ws.Range("K70").Value = 5938.857
$ws.Range("M70").Value = -5668.857
$ws.Range("H73").Value = 7156.9
$ws.Range("I73").Value = 5938.857
$ws.Range("K73").Value = 5938.857
$ws.Range("M73").Value = -5002.857
$ws.Range("H97").Value = 41667176
$ws.Range("I97").Value = 50000380
$ws.Range("J97").Value = 1157.75
$ws.Range("K97").Value = 50000380
$ws.Range("L97").Value = 1157.75
$ws.Range("M97").Value = -49999884
$ws.Range("N97").Value = -2149.75
$ws.Range("H109").Value = 32987.332
$ws.Range("J109").Value = 32987.332
$ws.Range("L109").Value = 32987.332
$ws.Range("N109").Value = -35067.332
$ws.Range("H113").Value = 7034.0527
$ws.Range("I113").Value = 7188.1177
$ws.Range("J113").Value = 5724.5
$ws.Range("K113").Value = 7188.1177
$ws.Range("L113").Value = 5724.5
$ws.Range("M113").Value = -5018.1177
$ws.Range("N113").Value = -10064.5
$ws.Range("H122").Value = 67452.07
$ws.Range("I122").Value = 102879.22
$ws.Range("K122").Value = 308637.66
$ws.Range("M122").Value = -306187.66
$ws.Range("H123").Value = 47949.285
$ws.Range("J123").Value = 47949.285
$ws.Range("L123").Value = 47949.285
$ws.Range("N123").Value = -52849.285
$ws.Range("H126").Value = 50983.316
$ws.Range("I126").Value = 67795.43
$ws.Range("J126").Value = 3909.4
$ws.Range("K126").Value = 203386.29
$ws.Range("L126").Value = 11728.2
$ws.Range("M126").Value = -200916.29
$ws.Range("N126").Value = -16668.2
$ws.Range("H132").Value = 2628.9512
$ws.Range("I132").Value = 2752
$ws.Range("J132").Value = 1490.75
$ws.Range("K132").Value = 8256
$ws.Range("L132").Value = 4472.25
$ws.Range("M132").Value = -5726
$ws.Range("N132").Value = -9532.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 20202.4
$ws.Range("I14").Value = 19001.334
$ws.Range("J14").Value = 22004
$ws.Range("K14").Value = 19001.334
$ws.Range("L14").Value = 22004
$ws.Range("M14").Value = -18829.334
$ws.Range("N14").Value = -22348
$ws.Range("H16").Value = 2103.4707
$ws.Range("I16").Value = 868.5
$ws.Range("J16").Value = 5067.4
$ws.Range("K16").Value = 868.5
$ws.Range("L16").Value = 5067.4
$ws.Range("M16").Value = -698.5
$ws.Range("N16").Value = -5407.4
$ws.Range("H40").Value = 19615.617
$ws.Range("I40").Value = 23046.12
$ws.Range("J40").Value = 10086.444
$ws.Range("K40").Value = 23046.12
$ws.Range("L40").Value = 10086.444
$ws.Range("M40").Value = -22910.12
$ws.Range("N40").Value = -10358.444
$ws.Range("H46").Value = 1686.8636
$ws.Range("J46").Value = 2188.4
$ws.Range("L46").Value = 2188.4
$ws.Range("N46").Value = -2564.4
$ws.Range("H93").Value = 1501.8096
$ws.Range("I93").Value = 1071.1538
$ws.Range("J93").Value = 2201.625
$ws.Range("K93").Value = 1071.1538
$ws.Range("L93").Value = 2201.625
$ws.Range("M93").Value = 176.8462
$ws.Range("N93").Value = -4697.625
$ws.Range("H94").Value = 77500
$ws.Range("J94").Value = 77500
$ws.Range("L94").Value = 77500
$ws.Range("N94").Value = -78852
$ws.Range("H100").Value = 1881.1538
$ws.Range("I100").Value = 1995
$ws.Range("J100").Value = 1625
$ws.Range("K100").Value = 1995
$ws.Range("L100").Value = 1625
$ws.Range("M100").Value = -1454
$ws.Range("N100").Value = -2707
$ws.Range("H122").Value = 3291.3713
$ws.Range("I122").Value = 3276.4119
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 9829.235700000001
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -7379.235700000001
$ws.Range("N122").Value = -16300
$ws.Range("H132").Value = 3518.261
$ws.Range("I132").Value = 2917.5
$ws.Range("J132").Value = 4891.4287
$ws.Range("K132").Value = 8752.5
$ws.Range("L132").Value = 14674.2861
$ws.Range("M132").Value = -6222.5
$ws.Range("N132").Value = -19734.2861
$ws.Range("H136").Value = 3869.6956
$ws.Range("I136").Value = 2181
$ws.Range("J136").Value = 9949
$ws.Range("K136").Value = 6543
$ws.Range("L136").Value = 29847
$ws.Range("M136").Value = -3993
$ws.Range("N136").Value = -34947
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 29133
$ws.Range("I54").Value = 29400
$ws.Range("K54").Value = 29400
$ws.Range("M54").Value = -28880
$ws.Range("H58").Value = 8904.333
$ws.Range("I58").Value = 8904.333
$ws.Range("K58").Value = 8904.333
$ws.Range("M58").Value = -8596.333
$ws.Range("H95").Value = 40664
$ws.Range("J95").Value = 40664
$ws.Range("L95").Value = 40664
$ws.Range("N95").Value = -46156
$ws.Range("H100").Value = 846.0741
$ws.Range("I100").Value = 797.0476
$ws.Range("K100").Value = 1594.0952
$ws.Range("M100").Value = -1053.0952
$ws.Range("H110").Value = 65333.332
$ws.Range("J110").Value = 65333.332
$ws.Range("L110").Value = 65333.332
$ws.Range("N110").Value = -73513.332
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H119").Value = 55000
$ws.Range("J119").Value = 55000
$ws.Range("L119").Value = 55000
$ws.Range("N119").Value = -64676
$ws.Range("H122").Value = 2523.1333
$ws.Range("I122").Value = 2649.7856
$ws.Range("K122").Value = 7949.3568
$ws.Range("M122").Value = -5499.3568
$ws.Range("H126").Value = 1986.75
$ws.Range("I126").Value = 1984.7
$ws.Range("J126").Value = 1997
$ws.Range("K126").Value = 5954.1
$ws.Range("L126").Value = 5991
$ws.Range("M126").Value = -3484.1
$ws.Range("N126").Value = -10931
$ws.Range("H132").Value = 4780.6665
$ws.Range("I132").Value = 3581.4546
$ws.Range("J132").Value = 5795.385
$ws.Range("K132").Value = 10744.3638
$ws.Range("L132").Value = 17386.155
$ws.Range("M132").Value = -8214.3638
$ws.Range("N132").Value = -22446.155
$ws.Range("H136").Value = 4336.1567
$ws.Range("I136").Value = 4079.0256
$ws.Range("K136").Value = 12237.0768
$ws.Range("M136").Value = -9687.076799999999
